# Atualização automática de PAROBE.xlsx
#
# - Renames "Paineis DARQ"            -> "PAINEIS DARQ"
# - Renames "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# - Removes the "Desarquivamentos Pendentes" worksheet entirely

$wb = $excel.ActiveWorkbook

# Rename the first sheet.
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"

# Rename the "Recolhimento x Eliminacao" sheet (now with accented caps).
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Drop the confirmation prompt Excel shows when deleting a sheet, then
# delete "Desarquivamentos Pendentes" outright.
$excel.DisplayAlerts = $false
$null = $wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()
$excel.DisplayAlerts = $true

# Leave the first sheet active/selected.
$wb.Worksheets.Item("PAINEIS DARQ").Activate()
